# Update TPM-derived statistics in Sema4d-Plxnb2 LR-pairs sheet
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("E2").Value = 2
$ws.Range("F2").Value = 0.6666666666666666
$ws.Range("G2").Value = 0.2115286666666667
$ws.Range("H2").Value = 0.634586
$ws.Range("I2").Value = 0.08153347995807345
$ws.Range("J2").Value = 0.08153347995807345
$ws.Range("K2").Value = 3
$ws.Range("L2").Value = 1
$ws.Range("M2").Value = 2.325008666666667
$ws.Range("N2").Value = 6.975026
$ws.Range("O2").Value = 0.05445297772988467
$ws.Range("P2").Value = 0.05445297772988466
$ws.Range("Q2").Value = 0.4918059832484445
$ws.Range("R2").Value = 4.426253849236
$ws.Range("S2").Value = 0.004439740768396972
$ws.Range("T2").Value = 0.00443974076839697

# Row 3
$ws.Range("E3").Value = 2
$ws.Range("F3").Value = 0.6666666666666666
$ws.Range("G3").Value = 0.2115286666666667
$ws.Range("H3").Value = 0.634586
$ws.Range("I3").Value = 0.08153347995807345
$ws.Range("J3").Value = 0.08153347995807345
$ws.Range("K3").Value = 3
$ws.Range("L3").Value = 1
$ws.Range("M3").Value = 19.33828433333333
$ws.Range("N3").Value = 58.014853
$ws.Range("O3").Value = 0.4529132218878514
$ws.Range("P3").Value = 0.4529132218878514
$ws.Range("Q3").Value = 4.090601500650889
$ws.Range("R3").Value = 36.815413505858
$ws.Range("S3").Value = 0.03692759109953961
$ws.Range("T3").Value = 0.0369275910995396

# Row 4
$ws.Range("E4").Value = 2
$ws.Range("F4").Value = 0.6666666666666666
$ws.Range("G4").Value = 0.2115286666666667
$ws.Range("H4").Value = 0.634586
$ws.Range("I4").Value = 0.08153347995807345
$ws.Range("J4").Value = 0.08153347995807345
$ws.Range("K4").Value = 3
$ws.Range("L4").Value = 1
$ws.Range("M4").Value = 21.03425566666667
$ws.Range("N4").Value = 63.102767
$ws.Range("O4").Value = 0.492633800382264
$ws.Range("P4").Value = 0.492633800382264
$ws.Range("Q4").Value = 4.449348055495778
$ws.Range("R4").Value = 40.044132499462
$ws.Range("S4").Value = 0.04016614809013688
$ws.Range("T4").Value = 0.04016614809013688

# Row 5
$ws.Range("E5").Value = 3
$ws.Range("F5").Value = 1
$ws.Range("G5").Value = 1.228643
$ws.Range("H5").Value = 3.685929
$ws.Range("I5").Value = 0.4735790235655714
$ws.Range("J5").Value = 0.4735790235655714
$ws.Range("K5").Value = 3
$ws.Range("L5").Value = 1
$ws.Range("M5").Value = 2.325008666666667
$ws.Range("N5").Value = 6.975026
$ws.Range("O5").Value = 0.05445297772988467
$ws.Range("P5").Value = 0.05445297772988466
$ws.Range("Q5").Value = 2.856605623239333
$ws.Range("R5").Value = 25.709450609154
$ws.Range("S5").Value = 0.02578778802355659
$ws.Range("T5").Value = 0.02578778802355658

# Row 6
$ws.Range("E6").Value = 3
$ws.Range("F6").Value = 1
$ws.Range("G6").Value = 1.228643
$ws.Range("H6").Value = 3.685929
$ws.Range("I6").Value = 0.4735790235655714
$ws.Range("J6").Value = 0.4735790235655714
$ws.Range("K6").Value = 3
$ws.Range("L6").Value = 1
$ws.Range("M6").Value = 19.33828433333333
$ws.Range("N6").Value = 58.014853
$ws.Range("O6").Value = 0.4529132218878514
$ws.Range("P6").Value = 0.4529132218878514
$ws.Range("Q6").Value = 23.75984767815967
$ws.Range("R6").Value = 213.838629103437
$ws.Range("S6").Value = 0.2144902013815857
$ws.Range("T6").Value = 0.2144902013815856

# Row 7
$ws.Range("E7").Value = 3
$ws.Range("F7").Value = 1
$ws.Range("G7").Value = 1.228643
$ws.Range("H7").Value = 3.685929
$ws.Range("I7").Value = 0.4735790235655714
$ws.Range("J7").Value = 0.4735790235655714
$ws.Range("K7").Value = 3
$ws.Range("L7").Value = 1
$ws.Range("M7").Value = 21.03425566666667
$ws.Range("N7").Value = 63.102767
$ws.Range("O7").Value = 0.492633800382264
$ws.Range("P7").Value = 0.492633800382264
$ws.Range("Q7").Value = 25.84359098506033
$ws.Range("R7").Value = 232.592318865543
$ws.Range("S7").Value = 0.2333010341604292
$ws.Range("T7").Value = 0.2333010341604292

# Row 8
$ws.Range("E8").Value = 3
$ws.Range("F8").Value = 1
$ws.Range("G8").Value = 1.154206333333333
$ws.Range("H8").Value = 3.462619
$ws.Range("I8").Value = 0.4448874964763552
$ws.Range("J8").Value = 0.4448874964763552
$ws.Range("K8").Value = 3
$ws.Range("L8").Value = 1
$ws.Range("M8").Value = 2.325008666666667
$ws.Range("N8").Value = 6.975026
$ws.Range("O8").Value = 0.05445297772988467
$ws.Range("P8").Value = 0.05445297772988466
$ws.Range("Q8").Value = 2.683539728121556
$ws.Range("R8").Value = 24.151857553094
$ws.Range("S8").Value = 0.02422544893793112
$ws.Range("T8").Value = 0.02422544893793111

# Row 9
$ws.Range("E9").Value = 3
$ws.Range("F9").Value = 1
$ws.Range("G9").Value = 1.154206333333333
$ws.Range("H9").Value = 3.462619
$ws.Range("I9").Value = 0.4448874964763552
$ws.Range("J9").Value = 0.4448874964763552
$ws.Range("K9").Value = 3
$ws.Range("L9").Value = 1
$ws.Range("M9").Value = 19.33828433333333
$ws.Range("N9").Value = 58.014853
$ws.Range("O9").Value = 0.4529132218878514
$ws.Range("P9").Value = 0.4529132218878514
$ws.Range("Q9").Value = 22.32037025333412
$ws.Range("R9").Value = 200.883332280007
$ws.Range("S9").Value = 0.2014954294067262
$ws.Range("T9").Value = 0.2014954294067262

# Row 10
$ws.Range("E10").Value = 3
$ws.Range("F10").Value = 1
$ws.Range("G10").Value = 1.154206333333333
$ws.Range("H10").Value = 3.462619
$ws.Range("I10").Value = 0.4448874964763552
$ws.Range("J10").Value = 0.4448874964763552
$ws.Range("K10").Value = 3
$ws.Range("L10").Value = 1
$ws.Range("M10").Value = 21.03425566666667
$ws.Range("N10").Value = 63.102767
$ws.Range("O10").Value = 0.492633800382264
$ws.Range("P10").Value = 0.492633800382264
$ws.Range("Q10").Value = 24.27787110741922
$ws.Range("R10").Value = 218.500839966773
$ws.Range("S10").Value = 0.2191666181316979
$ws.Range("T10").Value = 0.2191666181316979
